$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.412.35"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "1.643.62"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "299.05"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3785"
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3512"
$ws.Range("E8").Value = "  -3.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.86"
$ws.Range("E9").Value = "  -2.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08072"
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.210"
$ws.Range("E11").Value = "  -3.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.02"
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.348"
$ws.Range("E14").Value = "  -3.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.302"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001200"
$ws.Range("E16").Value = "  -3.39%  "
$ws.Range("D17").Value = "1.637.88"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.18"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06953"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.700"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.31"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.33"
$ws.Range("E23").Value = "  -3.87%  "
$ws.Range("D24").Value = "23.442.26"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.481"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.922"
$ws.Range("E26").Value = "  -4.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.84"
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.26"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.190"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.63"
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("D31").Value = "1.828.11"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.828"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.126"
$ws.Range("E33").Value = "  -6.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.49"
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9809"
$ws.Range("E35").Value = "  -8.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02699"
$ws.Range("E36").Value = "  -4.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08725"
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2426"
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.880"
$ws.Range("E39").Value = "  -4.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06804"
$ws.Range("E40").Value = "  -4.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.87"
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6847"
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.312"
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.57"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6343"
$ws.Range("E46").Value = "  -3.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.245"
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.905"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07706"
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.12"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.143"
$ws.Range("E51").Value = "  -4.33%  "
